$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 4: Title -> "test video", status -> "Upload"
$ws.Cells.Item(4, 3).Value = "test video"
$ws.Cells.Item(4, 7).Value = "Upload"

# Update row 11: Channel -> "Mushroom Toy Unboxing"
$ws.Cells.Item(11, 2).Value = "Mushroom Toy Unboxing"

# Remove the entire "URL" column (H) - this also drops the H1 header,
# all URL/time values in H2:H11, and shifts the dimension to A1:G11.
$ws.Columns.Item(8).Delete()
